$d = $word.ActiveDocument

# Update the date heading (first paragraph)
$p = $d.Paragraphs.Item(1)
$p.Range.Text = "2025-12-14 Sunday"

# Update each table cell value (row-major order, 20 rows x 5 cols)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "43+56="
$t.Cell(1, 2).Range.Text = "48+35="
$t.Cell(1, 3).Range.Text = "71-62="
$t.Cell(1, 4).Range.Text = "4+57="
$t.Cell(1, 5).Range.Text = "82-24="
$t.Cell(2, 1).Range.Text = "29+15="
$t.Cell(2, 2).Range.Text = "99-31="
$t.Cell(2, 3).Range.Text = "47-46="
$t.Cell(2, 4).Range.Text = "80-11="
$t.Cell(2, 5).Range.Text = "72-22="
$t.Cell(3, 1).Range.Text = "70-4="
$t.Cell(3, 2).Range.Text = "53+36="
$t.Cell(3, 3).Range.Text = "42+14="
$t.Cell(3, 4).Range.Text = "62+16="
$t.Cell(3, 5).Range.Text = "45-4="
$t.Cell(4, 1).Range.Text = "54-8="
$t.Cell(4, 2).Range.Text = "65-18="
$t.Cell(4, 3).Range.Text = "92-19="
$t.Cell(4, 4).Range.Text = "12+78="
$t.Cell(4, 5).Range.Text = "59-52="
$t.Cell(5, 1).Range.Text = "43+48="
$t.Cell(5, 2).Range.Text = "5-4="
$t.Cell(5, 3).Range.Text = "44-22="
$t.Cell(5, 4).Range.Text = "62-39="
$t.Cell(5, 5).Range.Text = "64-46="
$t.Cell(6, 1).Range.Text = "27+17="
$t.Cell(6, 2).Range.Text = "26+20="
$t.Cell(6, 3).Range.Text = "65+5="
$t.Cell(6, 4).Range.Text = "7+68="
$t.Cell(6, 5).Range.Text = "91-70="
$t.Cell(7, 1).Range.Text = "53+3="
$t.Cell(7, 2).Range.Text = "49-18="
$t.Cell(7, 3).Range.Text = "13+4="
$t.Cell(7, 4).Range.Text = "68-10="
$t.Cell(7, 5).Range.Text = "92-77="
$t.Cell(8, 1).Range.Text = "26+48="
$t.Cell(8, 2).Range.Text = "95-54="
$t.Cell(8, 3).Range.Text = "60-6="
$t.Cell(8, 4).Range.Text = "16+15="
$t.Cell(8, 5).Range.Text = "48+1="
$t.Cell(9, 1).Range.Text = "81-71="
$t.Cell(9, 2).Range.Text = "98-94="
$t.Cell(9, 3).Range.Text = "93-13="
$t.Cell(9, 4).Range.Text = "57-3="
$t.Cell(9, 5).Range.Text = "14+4="
$t.Cell(10, 1).Range.Text = "97-18="
$t.Cell(10, 2).Range.Text = "84-78="
$t.Cell(10, 3).Range.Text = "51-2="
$t.Cell(10, 4).Range.Text = "70-13="
$t.Cell(10, 5).Range.Text = "48+16="
$t.Cell(11, 1).Range.Text = "99-8="
$t.Cell(11, 2).Range.Text = "88-38="
$t.Cell(11, 3).Range.Text = "67+22="
$t.Cell(11, 4).Range.Text = "36-32="
$t.Cell(11, 5).Range.Text = "26-18="
$t.Cell(12, 1).Range.Text = "81-56="
$t.Cell(12, 2).Range.Text = "22+20="
$t.Cell(12, 3).Range.Text = "38+32="
$t.Cell(12, 4).Range.Text = "62-36="
$t.Cell(12, 5).Range.Text = "72+7="
$t.Cell(13, 1).Range.Text = "78+13="
$t.Cell(13, 2).Range.Text = "72-10="
$t.Cell(13, 3).Range.Text = "69-35="
$t.Cell(13, 4).Range.Text = "62+33="
$t.Cell(13, 5).Range.Text = "4+67="
$t.Cell(14, 1).Range.Text = "45+10="
$t.Cell(14, 2).Range.Text = "68-24="
$t.Cell(14, 3).Range.Text = "93-53="
$t.Cell(14, 4).Range.Text = "25+67="
$t.Cell(14, 5).Range.Text = "92-68="
$t.Cell(15, 1).Range.Text = "32+15="
$t.Cell(15, 2).Range.Text = "91-22="
$t.Cell(15, 3).Range.Text = "42+32="
$t.Cell(15, 4).Range.Text = "77-6="
$t.Cell(15, 5).Range.Text = "23+36="
$t.Cell(16, 1).Range.Text = "72+16="
$t.Cell(16, 2).Range.Text = "44+13="
$t.Cell(16, 3).Range.Text = "22+12="
$t.Cell(16, 4).Range.Text = "61-26="
$t.Cell(16, 5).Range.Text = "42+40="
$t.Cell(17, 1).Range.Text = "91-46="
$t.Cell(17, 2).Range.Text = "68-45="
$t.Cell(17, 3).Range.Text = "54+33="
$t.Cell(17, 4).Range.Text = "35+11="
$t.Cell(17, 5).Range.Text = "43-36="
$t.Cell(18, 1).Range.Text = "98-8="
$t.Cell(18, 2).Range.Text = "75+6="
$t.Cell(18, 3).Range.Text = "90-38="
$t.Cell(18, 4).Range.Text = "3+96="
$t.Cell(18, 5).Range.Text = "11+59="
$t.Cell(19, 1).Range.Text = "76-57="
$t.Cell(19, 2).Range.Text = "6+69="
$t.Cell(19, 3).Range.Text = "17+70="
$t.Cell(19, 4).Range.Text = "81+14="
$t.Cell(19, 5).Range.Text = "43-23="
$t.Cell(20, 1).Range.Text = "68-30="
$t.Cell(20, 2).Range.Text = "97-94="
$t.Cell(20, 3).Range.Text = "89-35="
$t.Cell(20, 4).Range.Text = "29-25="
$t.Cell(20, 5).Range.Text = "89-33="
